$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Majorelle Magdy'
$ws.Range("G3").Value = 'Dr. Asmaa Reda, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Amira Sobhy'
$ws.Range("G4").Value = 'Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Hend Mahmoud'
$ws.Range("G5").Value = 'Dr. Hanan Ragab, Dr. Nourhan Mahmoud, Dr. Nesma, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Mohammad El-Tanany'
$ws.Range("G6").Value = 'Dr. Asmaa Reda, Dr. Nahla Nagiub, Dr. Nourhan Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Amira Sobhy, Dr. Hend Mahmoud'
$ws.Range("G7").Value = 'Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Amira Sobhy, Dr. Hend Mahmoud'
$ws.Range("G8").Value = 'Dr. Asmaa Reda, Administrator, Dr. Manar Montaser, Dr. Shimaa Ahmad Mekki, Dr. Majorelle Magdy, Dr. Eman Tantawi'
$ws.Range("G9").Value = 'Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Manar Montaser, Dr. Majorelle Magdy, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Amira Sobhy, Dr. Hend Mahmoud'
$ws.Range("G10").Value = 'Dr. Alshimaa Atef, Dr. Rana Abo-Zaid, Dr. Shimaa Ahmad Mekki, Dr. Gehan Adel, Dr. Heba Mahmoud Ali, Dr. Sara Wael, Dr. Servinaz Sayed Mohammad'
$ws.Range("G11").Value = 'Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Hend Mahmoud'
$ws.Range("G12").Value = 'Administrator, Dr. Salma El-Gendy'
$ws.Range("G13").Value = 'Dr. Safa Hany, D Wessam Atef, Dr. Shimaa Ashraf, Dr. Mariam Nour El-Din, Dr. Omnia Mohammad'
$ws.Range("G14").Value = 'Dr. Safa Hany, Dr. Shimaa Ashraf'
$ws.Range("G17").Value = 'Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Marwa Mustafa, Dr. Eman M. Abo-Sakaya, Dr. Basma Hamed, Dr. Sarah Abdelmohsen, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Dina Adel, Dr. Nourhan Osama'
$ws.Range("G22").Value = 'Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed'
$ws.Range("G23").Value = 'Dr. Nourham Mostafa, Dr. Hana Amr'
$ws.Range("G24").Value = 'Dr. Wafaa Ebida, Dr. Marina Atef, Dr. Remon, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Ola Abd Al-Fattah, Dr. Yasmin, Dr. Youstina Magdy, Dr. Aya Emad, Dr. Maryam Ashraf, Dr. Monica'
$ws.Range("G25").Value = 'Dr. Remon, Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Youstina Magdy, Dr. Aya Emad, Dr. Marina Atef'
$ws.Range("G26").Value = 'Dr. Youstina Magdy, Dr. Gehad Salah'
$ws.Range("G27").Value = 'Dr. Wafaa Ebida, Dr. Remon, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Ola Abd Al-Fattah, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Eman Mohammad Al'
$ws.Range("G28").Value = 'Dr. Wafaa Ebida, Dr. Remon, Dr. Neveen Nashaat, Dr. Aya Hanafy, Dr. Nardine, Dr. Abdullah El-Agrody, Dr. Salma Hassan, Dr. Eman Samir Gabry'
$ws.Range("G29").Value = 'Dr. Remon, Dr. Neveen Nashaat, Dr. Naema Gomaa, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Monica'
$ws.Range("G30").Value = 'Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Hend Mahmoud'
$ws.Range("G31").Value = 'Dr. Asmaa Reda, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Amira Sobhy'
$ws.Range("G32").Value = 'Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Hend Mahmoud'
$ws.Range("G33").Value = 'Dr. Hanan Ragab, Dr. Nourhan Mahmoud, Dr. Nesma, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Mohammad El-Tanany'
$ws.Range("G34").Value = 'Dr. Asmaa Reda, Dr. Nahla Nagiub, Dr. Nourhan Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Amira Sobhy, Dr. Hend Mahmoud'
$ws.Range("G35").Value = 'Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Amira Sobhy, Dr. Hend Mahmoud'
$ws.Range("G36").Value = 'Dr. Asmaa Reda, Administrator, Dr. Manar Montaser, Dr. Shimaa Ahmad Mekki, Dr. Majorelle Magdy, Dr. Eman Tantawi'
$ws.Range("G37").Value = 'Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Manar Montaser, Dr. Majorelle Magdy, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Amira Sobhy, Dr. Hend Mahmoud'
$ws.Range("G38").Value = 'Dr. Alshimaa Atef, Dr. Rana Abo-Zaid, Dr. Shimaa Ahmad Mekki, Dr. Gehan Adel, Dr. Heba Mahmoud Ali, Dr. Sara Wael, Dr. Servinaz Sayed Mohammad'
$ws.Range("G39").Value = 'Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Hend Mahmoud'
$ws.Range("G40").Value = 'Administrator, Dr. Salma El-Gendy'
$ws.Range("G41").Value = 'Dr. Safa Hany, D Wessam Atef, Dr. Shimaa Ashraf, Dr. Mariam Nour El-Din, Dr. Omnia Mohammad'
$ws.Range("G42").Value = 'Dr. Safa Hany, Dr. Shimaa Ashraf'
$ws.Range("G45").Value = 'Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Marwa Mustafa, Dr. Eman M. Abo-Sakaya, Dr. Basma Hamed, Dr. Sarah Abdelmohsen, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Dina Adel, Dr. Nourhan Osama'
$ws.Range("G50").Value = 'Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed'
$ws.Range("G51").Value = 'Dr. Nourham Mostafa, Dr. Hana Amr'
$ws.Range("G52").Value = 'Dr. Wafaa Ebida, Dr. Marina Atef, Dr. Remon, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Ola Abd Al-Fattah, Dr. Yasmin, Dr. Youstina Magdy, Dr. Aya Emad, Dr. Maryam Ashraf, Dr. Monica'
$ws.Range("G53").Value = 'Dr. Remon, Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Youstina Magdy, Dr. Aya Emad, Dr. Marina Atef'
$ws.Range("G54").Value = 'Dr. Youstina Magdy, Dr. Gehad Salah'
$ws.Range("G55").Value = 'Dr. Wafaa Ebida, Dr. Remon, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Ola Abd Al-Fattah, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Eman Mohammad Al'
$ws.Range("G56").Value = 'Dr. Wafaa Ebida, Dr. Remon, Dr. Neveen Nashaat, Dr. Aya Hanafy, Dr. Nardine, Dr. Abdullah El-Agrody, Dr. Salma Hassan, Dr. Eman Samir Gabry'
$ws.Range("G57").Value = 'Dr. Remon, Dr. Neveen Nashaat, Dr. Naema Gomaa, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Monica'
